$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44242, "Primera", 45, 12000, 12000, 12000, "`$/caja 15 kilos granel", "Provincia de Limarí", 800, 15),
    @(3, 45001, "Primera", 50, 16000, 16000, 16000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1143, 14),
    @(4, 45050, "Especial", 56, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(5, 45050, "Primera", 50, 12000, 12000, 12000, "`$/caja 14 kilos granel", "Provincia de Limarí", 857, 14),
    @(6, 44323, "Primera", 60, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(7, 44627, "Primera", 56, 17000, 17000, 17000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1214, 14),
    @(8, 44592, "Primera", 54, 20000, 20000, 20000, "`$/caja 15 kilos empedrada", "Provincia de Limarí", 1333, 15),
    @(9, 44588, "Primera", 85, 19000, 20000, 19529, "`$/caja 14 kilos granel", "Provincia de Limarí", 1395, 14),
    @(10, 44320, "Primera", 45, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(11, 44252, "Primera", 60, 14000, 14000, 14000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1000, 14),
    @(12, 44314, "Primera", 56, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(13, 44322, "Primera", 50, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(14, 44278, "Primera", 45, 13000, 13000, 13000, "`$/caja 14 kilos empedrada", "Provincia del Elquí", 929, 14),
    @(15, 45042, "Especial", 50, 17000, 17000, 17000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1214, 14),
    @(16, 45042, "Primera", 50, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(17, 44239, "Primera", 70, 15000, 15000, 15000, "`$/caja 15 kilos granel", "Provincia de Limarí", 1000, 15),
    @(18, 44630, "Primera", 75, 15000, 15000, 15000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1071, 14),
    @(19, 44260, "Primera", 56, 13000, 13000, 13000, "`$/caja 14 kilos empedrada", "Provincia del Elquí", 929, 14),
    @(20, 44313, "Primera", 36, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(21, 44616, "Primera", 70, 14000, 14000, 14000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1000, 14),
    @(22, 44271, "Primera", 50, 12000, 12000, 12000, "`$/caja 14 kilos granel", "Provincia del Elquí", 857, 14),
    @(23, 44259, "Primera", 80, 12000, 12000, 12000, "`$/caja 15 kilos empedrada", "Provincia de Limarí", 800, 15),
    @(24, 44270, "Primera", 85, 12000, 12000, 12000, "`$/caja 14 kilos granel", "Provincia del Elquí", 857, 14),
    @(25, 45014, "Primera", 60, 15000, 15000, 15000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1071, 14),
    @(26, 45043, "Especial", 45, 17000, 17000, 17000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1214, 14),
    @(27, 45043, "Primera", 67, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(28, 45054, "Especial", 54, 16000, 16000, 16000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1143, 14),
    @(29, 45054, "Primera", 50, 14000, 14000, 14000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1000, 14),
    @(30, 44614, "Primera", 54, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(31, 44245, "Primera", 50, 15000, 15000, 15000, "`$/caja 15 kilos granel", "Provincia de Limarí", 1000, 15),
    @(32, 45040, "Especial", 65, 17000, 17000, 17000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1214, 14),
    @(33, 45040, "Primera", 60, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(34, 45006, "Primera", 40, 16000, 16000, 16000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1143, 14),
    @(35, 44316, "Primera", 48, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(36, 44238, "Primera", 60, 15000, 15000, 15000, "`$/caja 15 kilos granel", "Provincia de Limarí", 1000, 15),
    @(37, 44312, "Primera", 68, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(38, 45044, "Especial", 30, 16000, 16000, 16000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1143, 14),
    @(39, 45044, "Primera", 30, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
    @(40, 45015, "Primera", 56, 15000, 15000, 15000, "`$/caja 14 kilos empedrada", "Provincia de Limarí", 1071, 14),
    @(41, 44315, "Primera", 65, 14000, 14000, 14000, "`$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D Fecha
    $ws.Cells.Item($r, 12).Value = $row[2]   # L Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]   # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[5]   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[6]   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[7]   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $row[8]   # R Origen
    $ws.Cells.Item($r, 19).Value = $row[9]   # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[10]  # T Kg / unidad
}